# Update the dSF column (F) values for the severino_luis workbook.
# The diff shows that only column F (dSF) values changed for most data
# rows (rows 2-25), while rows 12 and 15 keep their original value of 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = -1
    3  = 2
    4  = 1
    5  = 0
    6  = 3
    7  = -4
    8  = 4
    9  = -1
    10 = -3
    11 = 5
    12 = 0
    13 = 2
    14 = -1
    15 = 0
    16 = 3
    17 = -3
    18 = -1
    19 = 3
    20 = 1
    21 = 1
    22 = -2
    23 = 7
    24 = -1
    25 = -4
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 6).Value = $newValues[$row]
}
